# Realestate Update resale numbers 2025-02-20 08:56
# Appends a new data row (row 79) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

# Columns A:D hold text that looks numeric/date-like (e.g. "2025-02-20",
# "08:56:21", "07"). Force text formatting first so Excel does not
# auto-convert them into date serials / numbers and drop the leading zero.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-20"
$ws.Cells.Item($row, 2).Value = "08:56:21"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "07"

# Restore the default "Normal" style so the cells don't keep a custom
# text number format applied to them (matches the original workbook
# where these cells carry no explicit style).
$textRange.Style = "Normal"

# Numeric columns E:T
$ws.Cells.Item($row, 5).Value = 129506
$ws.Cells.Item($row, 6).Value = 140536
$ws.Cells.Item($row, 7).Value = 171494
$ws.Cells.Item($row, 8).Value = 160021
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145919
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192575
$ws.Cells.Item($row, 14).Value = 115073
$ws.Cells.Item($row, 15).Value = 45815
$ws.Cells.Item($row, 16).Value = 29088
$ws.Cells.Item($row, 17).Value = 67186
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 46792
$ws.Cells.Item($row, 20).Value = -1
